# Applies the "Added a player mat" edit:
#   - Shifts the four existing shapes on slide 1 up-and-to-the-left by the
#     same vector (dx=-2422225 EMU, dy=-1114573 EMU) to make room for a new
#     player-mat graphic behind them.
#   - Appends a trailing tab character to the "ability (or else)" run.
#
# NOTE: PowerPoint's Shape.Left/Top/Width/Height are stored as single
# precision (32-bit) floats, so the point values below are chosen (to ~7
# significant decimal digits) so that, after the COM layer's internal
# float32 round-trip, they reproduce the exact target EMU offsets recorded
# in the OOXML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape #1: rotated background picture ("תמונה 5")
# a:off x="3772337" y="-244775" -> x="1350112" y="-1359348"
$s.Shapes.Item(1).Left = 106.30803299606299
$s.Shapes.Item(1).Top  = -107.03535079055118

# Shape #2: title textbox ("תיבת טקסט 3" - "Pokemon Munchkin Rules")
# a:off x="2635249" y="1514548" -> x="213024" y="399975"
$s.Shapes.Item(2).Left = 16.773543407086617
$s.Shapes.Item(2).Top  = 31.49409398818898

# Shape #3: bulleted rules textbox ("תיבת טקסט 6")
# a:off x="2905458" y="1976213" -> x="483233" y="861640"
$s.Shapes.Item(3).Left = 38.04984101968503
$s.Shapes.Item(3).Top  = 67.84566889133858

# Shape #4: small picture ("תמונה 8")
# a:off x="8591112" y="1212849" -> x="6168887" y="98276"
$s.Shapes.Item(4).Left = 485.73912055826776
$s.Shapes.Item(4).Top  = 7.738267716535433

# Append a trailing tab character to the "ability (or else)" run inside the
# bulleted rules textbox, preserving its existing run formatting.
$tr = $s.Shapes.Item(3).TextFrame.TextRange
$found = $tr.Find("ability (or else)")
if ($found -ne $null) {
    $found.Text = $found.Text + "`t"
}
